# 8d_cost_utilisation_intervention.xlsx — "small fixes + tornado diagram"
# Updates the per-row total-cost figures in column F across the five
# repeating 12-row blocks (rows 2-59) and moves the active selection down
# to G72 (scrolled towards the bottom of the sheet where new content was
# added).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 155.94 -> 128.1 (first row of each block + the "total" sub-row two below it)
$rows128 = @(2, 4, 14, 16, 26, 28, 38, 40, 50, 52)
foreach ($r in $rows128) {
    $ws.Cells.Item($r, 6).Value = 128.1
}

# 406.66999999999996 -> 378.83 (second row of each block)
$rows378 = @(3, 15, 27, 39, 51)
foreach ($r in $rows378) {
    $ws.Cells.Item($r, 6).Value = 378.83
}

# 474.66 -> 502.5 (the four-row "grouped" sub-block in the middle of each block)
$rows502 = @(8, 9, 10, 11, 20, 21, 22, 23, 32, 33, 34, 35, 44, 45, 46, 47, 56, 57, 58, 59)
foreach ($r in $rows502) {
    $ws.Cells.Item($r, 6).Value = 502.5
}

# Move the selection/scroll position towards the bottom of the sheet.
$ws.Range("G72").Select()
